$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" labels between row 16 and row 18 (2402 <-> 2407),
# and update "Valor Mora" (F) + "Salario Basico" (G) figures accordingly.

$ws.Range("E16").Value = "2407"
$ws.Range("F16").Value = 677761
$ws.Range("G16").Value = 24205754

$ws.Range("E17").Value = "2404"
$ws.Range("F17").Value = 632212
$ws.Range("G17").Value = 24205754

$ws.Range("E18").Value = "2402"
$ws.Range("F18").Value = 638622
$ws.Range("G18").Value = 24205754
